$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 554. This shifts the existing rows 554-622
# down to 555-623 (Excel copies formatting from the row above, which is
# what the source data already shows - e.g. the date-style on column D).
$ws.Rows.Item(554).Insert()

# Populate the newly inserted row 554 with the new record. Columns A, B,
# C, E, F, G, H, I, J, K, L, N, O, P, Q, R, S, T are identical to the
# (now shifted) row below it, so copy them from row 555; only D (Fecha)
# and M (Volumen) differ for the new record.
$ws.Range("A554:T554").Value = $ws.Range("A555:T555").Value2

$ws.Range("D554").Value = 45142
$ws.Range("M554").Value = 215
